$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 111702802
$ws.Range("Q2").Value = 516752
$ws.Range("R2").Value = 6574764
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = 111702796
$ws.Range("B3").Value = 90687
$ws.Range("E3").Value = 5964
$ws.Range("F3").Value = "Fjällig taggsvamp s.str."
$ws.Range("G3").Value = "Sarcodon imbricatus s.str."
$ws.Range("H3").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q3").Value = 516756
$ws.Range("R3").Value = 6574761
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# --- Row 4 ---
$ws.Range("A4").Value = 111702873
$ws.Range("B4").Value = 90332
$ws.Range("E4").Value = 4769
$ws.Range("F4").Value = "Svavelriska"
$ws.Range("G4").Value = "Lactarius scrobiculatus"
$ws.Range("H4").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q4").Value = 516761
$ws.Range("R4").Value = 6574773
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
